$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 483 (A_SERIES_5YAVERAGE): refine text to refer to the "Endjahr" (last year)
# instead of the "mittlere Jahr" (middle year).
$ws.Range("C483").Value = "Gleitender Fünfjahresdurchschnitt mit Bezug auf das Endjahr"
$ws.Range("D483").Value = "Moving five-year average shown for the last year"

# Re-key row 486 (A_SERIES_BEH) from K_SERIES to the new K_PRAEV category.
$ws.Range("B486").Value = "K_PRAEV"

# Re-key row 525 (A_SERIES_PREV) from K_ZUORDN to the new K_PRAEV category.
$ws.Range("B525").Value = "K_PRAEV"

# Remove row 526 (A_SERIES_PREVUNDBEH) entirely; all following rows shift up by one,
# which also reduces the sheet's used dimension from A1:D562 to A1:D561.
$ws.Rows(526).Delete()
